# Update for release-notes.md insert (f80ed2bb9e1dd81abc71d13817b8a44a756cee80)
# - bump metadata (version/status/date/contact)
# - swap the two "Mapping" columns on the Elements sheet (RIM Mapping <-> Spec metier)

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value  = "0.4.0-snapshot-1"
$meta.Range("B6").Value  = "draft"
$meta.Range("B8").Value  = "2024-05-23T12:16:26+00:00"
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)"

# --- Elements sheet ---------------------------------------------------------
# Columns AK (37) and AL (38) swap places: "Mapping: RIM Mapping" and
# "Mapping: Spécification métier vers l'extension ROR LocationResidentialCapacity"
# trade order, so every row's AK/AL value (and the column widths) are exchanged.
$el = $wb.Worksheets.Item("Elements")

$lastRow = 16
for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $el.Cells.Item($r, 37)
    $alCell = $el.Cells.Item($r, 38)
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    $akCell.Value = $alVal
    $alCell.Value = $akVal
}

$akWidth = $el.Columns.Item(37).ColumnWidth
$alWidth = $el.Columns.Item(38).ColumnWidth
$el.Columns.Item(37).ColumnWidth = $alWidth
$el.Columns.Item(38).ColumnWidth = $akWidth

Write-Host "edit applied"
